$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Remove the "Silly scenarios" scenario columns (E:J) from the Constants
# sheet, leaving only the "Base scenario" column (D). Unmerge the old
# "Silly scenarios" header cell and clear its text, then fully clear
# (contents + formatting) the rest of the now-unused scenario data so the
# formulas in column B keep referring to the (now empty) cells instead of
# being rewritten to #REF! as a structural column delete would cause.
$ws.Range("E2:I2").UnMerge()
$ws.Range("E2:I2").ClearContents()
$ws.Range("E3:J35").Clear()

$ws.Range("B19").Select()
